$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (5th column) - shifts Tag List etc. right
$ws.Range("E1").EntireColumn.Insert()

# Header for the new column
$ws.Range("E1").Value = "Generate YTD, Quarterly, Since Inception numbers"

# Fill the whole data range with "No" first (so "No" lands in the shared-strings
# table before "Yes" - matches the original authoring order), then flip the
# "Management Fees" row (row 2) to "Yes".
$ws.Range("E2").Value = "No"
$ws.Range("E3").Value = "No"
$ws.Range("E4").Value = "No"
$ws.Range("E5").Value = "No"
$ws.Range("E6").Value = "No"
$ws.Range("E7").Value = "No"
$ws.Range("E8").Value = "No"
$ws.Range("E9").Value = "No"
$ws.Range("E10").Value = "No"

# Row 2: Management Fees -> Yes
$ws.Range("E2").Value = "Yes"

# Fix column width for the newly inserted column E (target stored width 14.3125;
# the COM ColumnWidth setter quantizes to 1/7-character steps, so 13.57 is the
# closest reachable value -- it's also what the engine itself reports when
# reading back a column whose stored width is 14.3125, e.g. column D).
$ws.Range("E1").ColumnWidth = 13.57
# Column F naturally inherited the original column E's width (14.3125) from the
# column-insert shift, so it is left untouched.

# Update the selection to match the target state (activeCell=E3, sqref=E3)
$ws.Range("E3").Select()
